$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 12:38"

# Row 12 - Iran
$ws.Range("B12").Value = 96448
$ws.Range("C12").Value = 802
$ws.Range("D12").Value = 77350
$ws.Range("E12").Value = 12942
$ws.Range("F12").Value = 2787
$ws.Range("G12").Value = 65
$ws.Range("H12").Value = 6156

# Row 37 - Rumania
$ws.Range("B37").Value = 12732
$ws.Range("C37").Value = 165
$ws.Range("D37").Value = 4547
$ws.Range("E37").Value = 7430
$ws.Range("F37").Value = 265

# Row 54 - Finlandia
$ws.Range("E54").Value = 1956
$ws.Range("F54").Value = 52
$ws.Range("G54").Value = 2
$ws.Range("H54").Value = 220

# Row 55 - Marruecos
$ws.Range("B55").Value = 4687
$ws.Range("C55").Value = 118
$ws.Range("D55").Value = 1235
$ws.Range("E55").Value = 3280
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 172

# Row 90 - Hong Kong
$ws.Range("D90").Value = 864
$ws.Range("E90").Value = 172

# Row 99 - Libano
$ws.Range("B99").Value = 733
$ws.Range("C99").Value = 4
$ws.Range("E99").Value = 516
$ws.Range("F99").Value = 43
$ws.Range("G99").Value = 1
$ws.Range("H99").Value = 25

# Row 110 - San Marino
$ws.Range("D110").Value = 83
$ws.Range("E110").Value = 456
